# "final code for newtours"
#
# The third worksheet ("bookflightexcel") gains three new leading data
# columns (new shared strings "afreen" / "rahman" plus a repeated
# "Muslim"), which pushes every column from the old D onward three
# places to the right (D:P -> G:S). The workbook's active tab also
# moves from the second sheet ("flightfinderexcel") to this third
# sheet ("bookflightexcel").

$wb  = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("bookflightexcel")

# Insert three blank columns at D (old D:P shifts right to G:S).
$ws3.Range("D1:F1").EntireColumn.Insert() | Out-Null

# Populate the newly inserted cells.
$ws3.Range("D1").Value = "afreen"
$ws3.Range("E1").Value = "rahman"
$ws3.Range("F1").Value = "Muslim"

# bookflightexcel becomes the active sheet/tab, with I4 selected
# (was L4 before the column insert shifted things left by 3).
$ws3.Activate()
$ws3.Range("I4").Select() | Out-Null
